$d = $word.ActiveDocument

# --- Edit 1: merge "Best model: " + "Epoch 10:" runs into a single run ---
# Re-applying the same text via Find/Replace coalesces the two adjacent
# runs ("Best model: " and "Epoch 10:") into one run of "Best model: Epoch 10:".
$d.Content.Find.Execute("Best model: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Best model: ", 2)

# --- Edit 2: drop the Epoch 6-9 section (text + image) that followed the
# "Results are just 3-4% ..." paragraph, keeping the trailing bookmark ---

# Locate the "Results are just ..." paragraph and the final paragraph
# (which holds the _GoBack bookmark) by their known text.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Results are just 3-4%*") {
        $anchorIndex = $i
        break
    }
}

$lastIndex = $d.Paragraphs.Count

# Delete every whole paragraph strictly between the anchor paragraph and the
# very last paragraph (images, "Epoch 6:" .. "Epoch 9:" stats, etc.), leaving
# the last paragraph (which carries the "_GoBack" bookmark) intact for now.
if ($lastIndex -gt ($anchorIndex + 1)) {
    $startP = $d.Paragraphs($anchorIndex + 1)
    $endP = $d.Paragraphs($lastIndex - 1)
    $midRange = $d.Range($startP.Range.Start, $endP.Range.End)
    $midRange.Delete()
}

# Replace the remaining last paragraph's text ("Early stopping due to no
# improvement") with two plain spaces, keeping its run/bookmark structure.
$d.Content.Find.Execute("Early stopping due to no improvement", $true, $false, `
                         $false, $false, $false, $true, 1, $false, "  ", 2)

# Merge that (now two-space) paragraph into the anchor paragraph by deleting
# the paragraph mark that separates them.
$anchorP = $d.Paragraphs($anchorIndex)
$markRange = $d.Range($anchorP.Range.End - 1, $anchorP.Range.End)
$markRange.Delete()
